$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.840.26"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.888.17"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'0.7532"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.04%  "
$ws.Range("D6").Value = "'242.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.3123"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'25.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").Value = "'0.07119"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("D11").Value = "'0.08479"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.94%  "
$ws.Range("D12").Value = "'0.7601"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").Value = "1.886.82"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "'5.365"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "'6.133"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "29.932.63"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "'243.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("D20").Value = "'0.000007815"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "2.139.88"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'8.029"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").Value = "'9.373"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").Value = "'162.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("E30").Value = "  +3.67%  "
$ws.Range("D31").Value = "'1.534"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").Value = "'4.511"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").Value = "'4.139"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("D34").Value = "'0.05430"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("D35").Value = "'1.242"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "'0.7516"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "'1.002"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").Value = "'2.711"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").Value = "'0.01948"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("D40").Value = "'2.772"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "'0.4460"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'6.106"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("D43").Value = "1.093.90"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "'72.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").Value = "'0.8601"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'7.727"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("D48").Value = "'102.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "'1.860"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("D50").Value = "'3.050"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").Value = "2.037.67"
$ws.Range("E51").Value = "  +0.24%  "
